$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text (so numeric-looking strings like
# "238.26" or "6.480" are not coerced into Excel numbers, which would both
# change their type and could strip significant trailing zeros).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Column D (Price) updates
Set-TextValue "D2"  "238.26"
Set-TextValue "D3"  "21.76"
Set-TextValue "D4"  "5.385"
Set-TextValue "D5"  "0.05642"
Set-TextValue "D6"  "6.480"
Set-TextValue "D7"  "3.355"
Set-TextValue "D8"  "0.7958"
Set-TextValue "D9"  "1.024"
Set-TextValue "D11" "0.07322"
Set-TextValue "D12" "0.03172"
Set-TextValue "D13" "0.02975"
Set-TextValue "D14" "0.09250"
Set-TextValue "D15" "0.001674"
Set-TextValue "D16" "3.265"
Set-TextValue "D18" "0.0005748"
Set-TextValue "D19" "0.006249"
Set-TextValue "D20" "0.005085"
Set-TextValue "D23" "0.0004006"
Set-TextValue "D24" "3.890"
Set-TextValue "D26" "0.3298"
Set-TextValue "D40" "0.04094"
Set-TextValue "D41" "0.006913"
Set-TextValue "D44" "0.009422"
Set-TextValue "D45" "0.00005444"
Set-TextValue "D47" "0.6762"
Set-TextValue "D48" "0.03796"

# Column E (Volume(1h)) updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
